$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet contained three duplicated rows (10/11, 12/13, 32/33 were each
# identical pairs). The update removes the extra (later) copy of each pair,
# which shifts every following row up by one. Delete from the bottom up so
# earlier row numbers stay stable while we work.
$ws.Rows.Item(33).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(10).Delete()

# Column A is a plain running index (0, 1, 2, ...) tied to row position, not
# to the row's content, so renumber it after the shift.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
